$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.009.32"
$ws.Range("E2").Value = "  -2.06%  "

$ws.Range("D3").Value = "3.061.79"
$ws.Range("E3").Value = "  -2.33%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.45%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "3.063.62"
$ws.Range("E8").Value = "  -2.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.466"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.05%  "

$ws.Range("E11").Value = "  -3.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.400"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.84%  "

$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").Value = "3.589.90"
$ws.Range("E14").Value = "  -2.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.02"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.72%  "

$ws.Range("E16").Value = "  -3.75%  "

$ws.Range("D17").Value = "57.071.31"
$ws.Range("E17").Value = "  -2.00%  "

$ws.Range("D18").Value = "3.065.55"
$ws.Range("E18").Value = "  -2.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "346.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.495"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.64%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.164"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.19%  "

$ws.Range("D28").Value = "0.0₃0851"
$ws.Range("E28").Value = "  -8.87%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0652"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.23%  "

$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("D44").Value = "2.396.77"
$ws.Range("E44").Value = "  +5.93%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("D47").Value = "3.105.27"
$ws.Range("E47").Value = "  -2.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0260"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.35%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.927"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.56%  "
